# Auto-generated edit script applying the Aegis_Profits.xlsx diff
# Updates computed market-price / profit columns (H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1332.4445
$ws.Range("I19").Value = 1175
$ws.Range("J19").Value = 1377.4286
$ws.Range("K19").Value = 1175
$ws.Range("L19").Value = 1377.4286
$ws.Range("M19").Value = -1000
$ws.Range("N19").Value = -1727.4286

$ws.Range("H86").Value = 16988
$ws.Range("I86").Value = 16220.8
$ws.Range("K86").Value = 16220.8
$ws.Range("M86").Value = -15097.8

$ws.Range("H89").Value = 16988
$ws.Range("I89").Value = 16220.8
$ws.Range("K89").Value = 81104
$ws.Range("M89").Value = -75488

$ws.Range("H129").Value = 793.3333
$ws.Range("J129").Value = 881.25
$ws.Range("L129").Value = 2643.75
$ws.Range("N129").Value = -12643.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 15000
$ws.Range("J51").Value = 15000
$ws.Range("L51").Value = 15000
$ws.Range("N51").Value = -16512

$ws.Range("H61").Value = 2732.4285
$ws.Range("I61").Value = 1196.8
$ws.Range("K61").Value = 1196.8
$ws.Range("M61").Value = -984.8

$ws.Range("H74").Value = 2338.1155
$ws.Range("I74").Value = 1272
$ws.Range("J74").Value = 5231.857
$ws.Range("K74").Value = 1272
$ws.Range("L74").Value = 5231.857
$ws.Range("M74").Value = -398
$ws.Range("N74").Value = -6979.857

$ws.Range("H77").Value = 2338.1155
$ws.Range("I77").Value = 1272
$ws.Range("J77").Value = 5231.857
$ws.Range("K77").Value = 6360
$ws.Range("L77").Value = 26159.285
$ws.Range("M77").Value = -1992
$ws.Range("N77").Value = -34895.285

$ws.Range("H136").Value = 2732.4285
$ws.Range("I136").Value = 1196.8
$ws.Range("K136").Value = 3590.4
$ws.Range("M136").Value = -1040.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 27002.719
$ws.Range("I20").Value = 38516.965
$ws.Range("J20").Value = 1095.6666
$ws.Range("K20").Value = 38516.965
$ws.Range("L20").Value = 1095.6666
$ws.Range("M20").Value = -38269.965
$ws.Range("N20").Value = -1589.6666

$ws.Range("H36").Value = 2018.5
$ws.Range("I36").Value = 2018.5
$ws.Range("K36").Value = 2018.5
$ws.Range("M36").Value = -1484.5

$ws.Range("H37").Value = 581.25
$ws.Range("I37").Value = 581.25
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 581.25
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -444.25
$ws.Range("N37").ClearContents()

$ws.Range("H134").Value = 1876.4839
$ws.Range("I134").Value = 2023.76
$ws.Range("J134").Value = 1262.8334
$ws.Range("K134").Value = 6071.28
$ws.Range("L134").Value = 3788.5002
$ws.Range("M134").Value = -3536.28
$ws.Range("N134").Value = -8858.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16778.031
$ws.Range("I31").Value = 22729.262
$ws.Range("J31").Value = 3090.2
$ws.Range("K31").Value = 22729.262
$ws.Range("L31").Value = 3090.2
$ws.Range("M31").Value = -22434.262
$ws.Range("N31").Value = -3680.2

$ws.Range("H34").Value = 16778.031
$ws.Range("I34").Value = 22729.262
$ws.Range("J34").Value = 3090.2
$ws.Range("K34").Value = 22729.262
$ws.Range("L34").Value = 3090.2
$ws.Range("M34").Value = -22527.262
$ws.Range("N34").Value = -3494.2

$ws.Range("H70").Value = 12749.75
$ws.Range("J70").Value = 12749.75
$ws.Range("L70").Value = 12749.75
$ws.Range("N70").Value = -13379.75

$ws.Range("H73").Value = 12749.75
$ws.Range("J73").Value = 12749.75
$ws.Range("L73").Value = 12749.75
$ws.Range("N73").Value = -14933.75

$ws.Range("H122").Value = 3349.3333
$ws.Range("I122").Value = 3095.9092
$ws.Range("J122").Value = 3856.182
$ws.Range("K122").Value = 9287.7276
$ws.Range("L122").Value = 11568.546
$ws.Range("M122").Value = -6837.7276
$ws.Range("N122").Value = -16468.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3832
$ws.Range("I3").Value = 3790
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 11370
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = -11258
$ws.Range("N3").Value = -12224

$ws.Range("H58").Value = 1148.75
$ws.Range("J58").Value = 1148.75
$ws.Range("L58").Value = 3446.25
$ws.Range("N58").Value = -3702.25

$ws.Range("H80").Value = 10395.333
$ws.Range("I80").Value = 1700
$ws.Range("J80").Value = 11185.818
$ws.Range("K80").Value = 5100
$ws.Range("L80").Value = 33557.454
$ws.Range("M80").Value = -4164
$ws.Range("N80").Value = -35429.454

$ws.Range("H83").Value = 10395.333
$ws.Range("I83").Value = 1700
$ws.Range("J83").Value = 11185.818
$ws.Range("K83").Value = 15300
$ws.Range("L83").Value = 100672.362
$ws.Range("M83").Value = -10620
$ws.Range("N83").Value = -110032.362

$ws.Range("H107").Value = 695722.1
$ws.Range("I107").Value = 1165
$ws.Range("J107").Value = 869361.4399999999
$ws.Range("K107").Value = 3495
$ws.Range("L107").Value = 2608084.32
$ws.Range("M107").Value = -1575
$ws.Range("N107").Value = -2611924.32

$ws.Range("H131").Value = 11223.371
$ws.Range("I131").Value = 880
$ws.Range("J131").Value = 11441.126
$ws.Range("K131").Value = 2640
$ws.Range("L131").Value = 34323.378
$ws.Range("M131").Value = 2400
$ws.Range("N131").Value = -44403.378

$ws.Range("H140").Value = 1842.6
$ws.Range("I140").Value = 1376.5555
$ws.Range("K140").Value = 4129.666499999999
$ws.Range("M140").Value = 1050.333500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 10366.143
$ws.Range("I52").Value = 5030
$ws.Range("J52").Value = 11255.5
$ws.Range("K52").Value = 5030
$ws.Range("L52").Value = 11255.5
$ws.Range("M52").Value = -4771
$ws.Range("N52").Value = -11773.5

$ws.Range("H102").Value = 208779
$ws.Range("I102").Value = 1959.64
$ws.Range("J102").Value = 1501400
$ws.Range("K102").Value = 1959.64
$ws.Range("L102").Value = 1501400
$ws.Range("M102").Value = -337.6400000000001
$ws.Range("N102").Value = -1504644

$ws.Range("H132").Value = 2201.4783
$ws.Range("I132").Value = 2139.4546
$ws.Range("J132").Value = 2258.3333
$ws.Range("K132").Value = 6418.3638
$ws.Range("L132").Value = 6774.999899999999
$ws.Range("M132").Value = -3888.3638
$ws.Range("N132").Value = -11834.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1875.25
$ws.Range("I136").Value = 1750.4445
$ws.Range("J136").Value = 2035.7142
$ws.Range("K136").Value = 5251.333500000001
$ws.Range("L136").Value = 6107.142599999999
$ws.Range("M136").Value = -2701.333500000001
$ws.Range("N136").Value = -11207.1426

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 6930.2856
$ws.Range("J49").Value = 6930.2856
$ws.Range("L49").Value = 6930.2856
$ws.Range("N49").Value = -7390.2856

$ws.Range("H54").Value = 6923.1333
$ws.Range("J54").Value = 6912.643
$ws.Range("L54").Value = 6912.643
$ws.Range("N54").Value = -7952.643

$ws.Range("H126").Value = 1145.125
$ws.Range("I126").Value = 1124.65
$ws.Range("J126").Value = 1247.5
$ws.Range("K126").Value = 3373.95
$ws.Range("L126").Value = 3742.5
$ws.Range("M126").Value = -903.9500000000003
$ws.Range("N126").Value = -8682.5

$ws.Range("H136").Value = 1528.359
$ws.Range("I136").Value = 414.4375
$ws.Range("J136").Value = 2303.261
$ws.Range("K136").Value = 1243.3125
$ws.Range("L136").Value = 6909.782999999999
$ws.Range("M136").Value = 1306.6875
$ws.Range("N136").Value = -12009.783

Write-Host "Applied 191 cell updates and 1 clear"